$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values: shrink the leading block to 16/20/16/20
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) tweaks
$ws.Range("B2").Value = 262.85707341393442
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 260.71443222478274
$ws.Range("E2").Value = 294.19274825932462

# Row 3 (STR) tweaks
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 292.59556661028108
$ws.Range("D3").Value = 256.34900221840485
$ws.Range("E3").Value = 301.36691324747699

# Selection now only covers the touched block
$ws.Range("B1:E3").Select()
